# New crime data collected — weekly CompStat update (022 Pct).
#
# Updates the "Volume / Number" and "Report Covering the Week ... Through ..."
# banner text, plus the weekly crime-complaint figures for rows 16-28
# (Murder .. Misd. Assault), including the TOTAL row 21.
#
# Some cells in this sheet alternate between a literal numeric value and a
# text placeholder ("0" or "***.*", used when a percent-change is undefined,
# e.g. division by zero). Plainly assigning a numeric-looking string via
# .Value lets Excel's smart-parsing turn it back into a number, so for the
# cells that must become TEXT we first force NumberFormat="@" before writing
# the string, then use Copy/PasteSpecial(xlPasteFormats) from a same-shaped
# reference cell to restore the real (General) display format without
# disturbing the stored text type. The reverse (text -> number) just needs
# the numeric-style format pasted in before the numeric value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-TextCell($addr, $formatSourceAddr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($formatSourceAddr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}

function Set-NumberCell($addr, $formatSourceAddr, $number) {
    $ws.Range($formatSourceAddr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).Value = $number
}

# ---- Masthead: Volume/Number and report week banner ----
$ws.Range("A8").Value = "Volume 31   Number  11"
$ws.Range("C9").Value = "Report Covering the Week  3/11/2024  Through  3/17/2024"

# ---- Row 16 (Murder) ----
Set-TextCell "C16" "C14" "0"
$ws.Range("F16").Value = 6
Set-TextCell "G16" "C14" "0"
Set-TextCell "H16" "C14" "***.*"
$ws.Range("L16").Value = 150

# ---- Row 17 (Rape) ----
Set-TextCell "C17" "C14" "0"
$ws.Range("E17").Value = -100
Set-TextCell "F17" "C14" "0"
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = -100
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 2
$ws.Range("K17").Value = -50
$ws.Range("L17").Value = -50
$ws.Range("N17").Value = -87.5

# ---- Row 18 (Robbery) ----
Set-NumberCell "D18" "I15" 1
Set-NumberCell "E18" "L15" -100
Set-NumberCell "G18" "I15" 1
Set-NumberCell "H18" "L15" -100
Set-NumberCell "J18" "I15" 1
Set-NumberCell "K18" "L15" -100

# ---- Row 19 (Fel. Assault) ----
$ws.Range("C19").Value = 3
$ws.Range("F19").Value = 4
$ws.Range("H19").Value = 300
$ws.Range("I19").Value = 6
$ws.Range("K19").Value = 50
$ws.Range("L19").Value = 100
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = -40

# ---- Row 21 (TOTAL) ----
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 50
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 150
$ws.Range("I21").Value = 18
$ws.Range("J21").Value = 9
$ws.Range("K21").Value = 100
$ws.Range("L21").Value = 63.636363636363
$ws.Range("M21").Value = 80
$ws.Range("N21").Value = -57.142857142857

# ---- Row 24 (Gr. Larceny) ----
Set-NumberCell "C24" "I15" 1
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 4
$ws.Range("I24").Value = 2
$ws.Range("J24").Value = 9
$ws.Range("K24").Value = -77.777777777777
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 100

# ---- Row 26 (Transit) ----
Set-NumberCell "C26" "I15" 1
Set-NumberCell "D26" "I15" 1
Set-NumberCell "E26" "L15" 0
$ws.Range("I26").Value = 4
$ws.Range("J26").Value = 4
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -42.857142857142

# ---- Row 28 (Petit Larceny) ----
Set-TextCell "C28" "C14" "0"

$ws.Application.CutCopyMode = $false
